$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update total_registros for row 2 (NEYRA PEREIRA MONICA AGNES ALEXANDRA)
$ws.Range("B2").Value = 145

# Rows 4, 6-11 are re-sorted/updated to keep the list in descending order
# by total_registros, with one name swapped in (TENE TRABUCCO GIAN PIERRE
# moved from row 11 to row 4) and updated counts for the rest.
$ws.Range("A4").Value = "TENE TRABUCCO GIAN PIERRE"

$ws.Range("A6").Value = "TUANAMA PIZANGO ELIZABETH"
$ws.Range("B6").Value = 113

$ws.Range("A7").Value = "SILVA ALVARADO EVELYN DE JESUS"
$ws.Range("B7").Value = 110

$ws.Range("A8").Value = "MORETO ESPINOZA CRISTIAN ESTEBAN"
$ws.Range("B8").Value = 107

$ws.Range("A9").Value = "SANCARRANCO SANCHEZ DE CRUZ GISSELA SHANI"
$ws.Range("B9").Value = 104

$ws.Range("A10").Value = "CHIROQUE YARLEQUE BETTY ELIZABETH"
$ws.Range("B10").Value = 102

$ws.Range("A11").Value = "HERRERA JUAN MANUEL"
$ws.Range("B11").Value = 102
